# DIV-4872 - Add UserRole column to CaseTypeTab
#
# Inserts a new "UserRole" column between the existing "TabFieldDisplayOrder"
# (col I) and "FieldShowCondition" (col J, was I) columns on the CaseTypeTab
# sheet, then makes CaseTypeTab the active/selected sheet (it was previously
# on ComplexTypes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CaseTypeTab")

# Shift columns I:K right to J:L, creating a blank column I and carrying
# over the existing cell styles from the old column I (F2:I2 / H3 share the
# same formatting the new column should inherit).
$ws.Columns("I:I").Insert()

# Populate the new column's header (row 2) and field name (row 3).
$ws.Range("I2").Value = "MaxLength: 100. No entry for role means no role restriction for that tab. Enter role on a single row per tab"
$ws.Range("I3").Value = "UserRole"

# Make CaseTypeTab the active sheet with I4 selected (previously ComplexTypes
# was the active/selected sheet).
$ws.Activate()
$ws.Range("I4").Select()
